$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-08-19 Monday" "2024-08-20 Tuesday"

Replace-Text "940×6=5640" "233×2=466"
Replace-Text "975×6=5850" "436×3=1308"
Replace-Text "141×8=1128" "921×9=8289"
Replace-Text "770×5=3850" "207×4=828"
Replace-Text "267×6=1602" "988×7=6916"
Replace-Text "682×6=4092" "150×5=750"
Replace-Text "934×9=8406" "638×6=3828"
Replace-Text "575×6=3450" "126×8=1008"
Replace-Text "103×3=309" "858×5=4290"
Replace-Text "788×6=4728" "462×9=4158"
Replace-Text "369×9=3321" "341×4=1364"
Replace-Text "822×6=4932" "103×6=618"
Replace-Text "727×6=4362" "868×5=4340"
Replace-Text "463×5=2315" "301×4=1204"
Replace-Text "347×7=2429" "233×8=1864"
Replace-Text "769×9=6921" "293×9=2637"
Replace-Text "266×9=2394" "758×5=3790"
Replace-Text "336×9=3024" "244×6=1464"
Replace-Text "988×3=2964" "179×5=895"
Replace-Text "355×4=1420" "152×6=912"
Replace-Text "383×3=1149" "616×9=5544"
Replace-Text "320×4=1280" "720×7=5040"
Replace-Text "720×8=5760" "621×7=4347"
Replace-Text "169×3=507" "817×2=1634"
Replace-Text "429×4=1716" "198×4=792"
